$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ewewereer"
$ws.Range("B1").Value = 45
$ws.Range("A2").Value = "qasdff"
$ws.Range("B2").Value = 5
$ws.Range("A3").Value = "aasfdf"
$ws.Range("B3").Value = 98

$ws.Range("B1").Select() | Out-Null
